$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 455
$ws.Range("F5").Value = 321
$ws.Range("G5").Value = 85
$ws.Range("F7").Value = 2120
$ws.Range("F8").Value = 47
$ws.Range("F9").Value = 48
$ws.Range("F10").Value = 1623
$ws.Range("F11").Value = 1623
$ws.Range("F12").Value = 1354
$ws.Range("F14").Value = 1404
$ws.Range("F15").Value = 19
$ws.Range("F17").Value = 559
$ws.Range("G17").Value = 149
$ws.Range("F18").Value = 153
$ws.Range("F20").Value = 7199
$ws.Range("F21").Value = 7865
$ws.Range("F22").Value = 47
$ws.Range("F28").Value = 217
$ws.Range("F35").Value = 1424
$ws.Range("F40").Value = 10
$ws.Range("F43").Value = 1360
$ws.Range("F46").Value = 191
$ws.Range("F48").Value = 167
$ws.Range("F49").Value = 159

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2616

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 455
$ws.Range("F7").Value = 321
$ws.Range("G7").Value = 85
$ws.Range("F10").Value = 2120
$ws.Range("F11").Value = 47
$ws.Range("F12").Value = 48
$ws.Range("F13").Value = 1623
$ws.Range("F14").Value = 1623
$ws.Range("F16").Value = 1354
$ws.Range("F18").Value = 1404
$ws.Range("F19").Value = 19
$ws.Range("F20").Value = 559
$ws.Range("G20").Value = 149
$ws.Range("F22").Value = 7199
$ws.Range("F23").Value = 7865
$ws.Range("F24").Value = 47
$ws.Range("F31").Value = 1424
$ws.Range("F44").Value = 191
$ws.Range("F46").Value = 167
$ws.Range("F47").Value = 159
